# [IMP] Add partner detail report
# Applies the template changes:
#  - fix typo "Activie" -> "Active" (column AD header)
#  - add 5 new report columns (AG:AK) with a new bold/green header style
#  - widen a few of the new columns
#  - refresh the view (zoom/scroll/selection)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- fix existing typo -------------------------------------------------
$ws.Range("AD5").Value = "Active"

# ---- new header cells ---------------------------------------------------
$ws.Range("AG5").Value = "Create Date"
$ws.Range("AH5").Value = "Create By"
$ws.Range("AI5").Value = "Update Date"
$ws.Range("AJ5").Value = "Update By"
$ws.Range("AK5").Value = "Active/Inactive"

# ---- style the first new header cell (bold, green fill, wrapped, centered)
$first = $ws.Range("AG5")
$first.Font.Name = "Arial"
$first.Font.Size = 10
$first.Font.Bold = $true
$first.Interior.Color = 11460802
$first.Interior.PatternColor = 12632256
$first.HorizontalAlignment = -4108
$first.VerticalAlignment = -4107
$first.WrapText = $true

# ---- copy that style onto the rest of the new header cells --------------
$first.Copy()
$ws.Range("AH5:AK5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---- widen the new columns (closest achievable given 1/6 char rounding) -
$ws.Columns("AG").ColumnWidth = 16.666666666666668
$ws.Columns("AH").ColumnWidth = 20.5
$ws.Columns("AJ").ColumnWidth = 19.666666666666668

# ---- refresh view: zoom out, scroll right, change active selection ------
$win = $excel.ActiveWindow
$win.Zoom = 100
$win.ScrollColumn = 23
$win.ScrollRow = 1
$ws.Range("AD11").Select()

Write-Host "done"
